$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 80.60050966666667
$ws.Cells.Item(2, 8).Value = 241.801529
$ws.Cells.Item(2, 9).Value = 0.17420496858261
$ws.Cells.Item(2, 10).Value = 0.17420496858261
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 0.8109183333333333
$ws.Cells.Item(2, 14).Value = 2.432755
$ws.Cells.Item(2, 15).Value = 0.2370884169621149
$ws.Cells.Item(2, 16).Value = 0.2370884169621149
$ws.Cells.Item(2, 17).Value = 65.36043096471056
$ws.Cells.Item(2, 18).Value = 588.243878682395
$ws.Cells.Item(2, 19).Value = 0.04130198022818597
$ws.Cells.Item(2, 20).Value = 0.04130198022818597
$ws.Cells.Item(3, 7).Value = 80.60050966666667
$ws.Cells.Item(3, 8).Value = 241.801529
$ws.Cells.Item(3, 9).Value = 0.17420496858261
$ws.Cells.Item(3, 10).Value = 0.17420496858261
$ws.Cells.Item(3, 15).Value = 0.3378801459239538
$ws.Cells.Item(3, 16).Value = 0.3378801459239539
$ws.Cells.Item(3, 17).Value = 93.14665066719722
$ws.Cells.Item(3, 18).Value = 838.319856004775
$ws.Cells.Item(3, 19).Value = 0.05886040020537005
$ws.Cells.Item(3, 20).Value = 0.05886040020537006
$ws.Cells.Item(4, 7).Value = 80.60050966666667
$ws.Cells.Item(4, 8).Value = 241.801529
$ws.Cells.Item(4, 9).Value = 0.17420496858261
$ws.Cells.Item(4, 10).Value = 0.17420496858261
$ws.Cells.Item(4, 13).Value = 1.434534666666667
$ws.Cells.Item(4, 14).Value = 4.303604
$ws.Cells.Item(4, 15).Value = 0.4194152964814894
$ws.Cells.Item(4, 16).Value = 0.4194152964814894
$ws.Cells.Item(4, 17).Value = 115.6242252678351
$ws.Cells.Item(4, 18).Value = 1040.618027410516
$ws.Cells.Item(4, 19).Value = 0.07306422854662391
$ws.Cells.Item(4, 20).Value = 0.07306422854662391
$ws.Cells.Item(5, 7).Value = 80.60050966666667
$ws.Cells.Item(5, 8).Value = 241.801529
$ws.Cells.Item(5, 9).Value = 0.17420496858261
$ws.Cells.Item(5, 10).Value = 0.17420496858261
$ws.Cells.Item(5, 11).Value = 1
$ws.Cells.Item(5, 12).Value = 0.3333333333333333
$ws.Cells.Item(5, 13).Value = 0.019209
$ws.Cells.Item(5, 14).Value = 0.057627
$ws.Cells.Item(5, 15).Value = 0.005616140632441737
$ws.Cells.Item(5, 16).Value = 0.005616140632441737
$ws.Cells.Item(5, 17).Value = 1.548255190187
$ws.Cells.Item(5, 18).Value = 13.934296711683
$ws.Cells.Item(5, 19).Value = 0.0009783596024300321
$ws.Cells.Item(5, 20).Value = 0.0009783596024300321
$ws.Cells.Item(6, 7).Value = 93.34790299999999
$ws.Cells.Item(6, 9).Value = 0.2017563980255169
$ws.Cells.Item(6, 10).Value = 0.2017563980255169
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 0.8109183333333333
$ws.Cells.Item(6, 14).Value = 2.432755
$ws.Cells.Item(6, 15).Value = 0.2370884169621149
$ws.Cells.Item(6, 16).Value = 0.2370884169621149
$ws.Cells.Item(6, 17).Value = 75.69752592092165
$ws.Cells.Item(6, 18).Value = 681.2777332882949
$ws.Cells.Item(6, 19).Value = 0.04783410501984816
$ws.Cells.Item(6, 20).Value = 0.04783410501984817
$ws.Cells.Item(7, 7).Value = 93.34790299999999
$ws.Cells.Item(7, 9).Value = 0.2017563980255169
$ws.Cells.Item(7, 10).Value = 0.2017563980255169
$ws.Cells.Item(7, 15).Value = 0.3378801459239538
$ws.Cells.Item(7, 16).Value = 0.3378801459239539
$ws.Cells.Item(7, 17).Value = 107.8782820011416
$ws.Cells.Item(7, 18).Value = 970.9045380102748
$ws.Cells.Item(7, 19).Value = 0.06816948120595295
$ws.Cells.Item(7, 20).Value = 0.06816948120595297
$ws.Cells.Item(8, 7).Value = 93.34790299999999
$ws.Cells.Item(8, 9).Value = 0.2017563980255169
$ws.Cells.Item(8, 10).Value = 0.2017563980255169
$ws.Cells.Item(8, 13).Value = 1.434534666666667
$ws.Cells.Item(8, 14).Value = 4.303604
$ws.Cells.Item(8, 15).Value = 0.4194152964814894
$ws.Cells.Item(8, 16).Value = 0.4194152964814894
$ws.Cells.Item(8, 17).Value = 133.9108029141373
$ws.Cells.Item(8, 18).Value = 1205.197226227236
$ws.Cells.Item(8, 19).Value = 0.08461971949490954
$ws.Cells.Item(8, 20).Value = 0.08461971949490954
$ws.Cells.Item(9, 7).Value = 93.34790299999999
$ws.Cells.Item(9, 9).Value = 0.2017563980255169
$ws.Cells.Item(9, 10).Value = 0.2017563980255169
$ws.Cells.Item(9, 11).Value = 1
$ws.Cells.Item(9, 12).Value = 0.3333333333333333
$ws.Cells.Item(9, 13).Value = 0.019209
$ws.Cells.Item(9, 14).Value = 0.057627
$ws.Cells.Item(9, 15).Value = 0.005616140632441737
$ws.Cells.Item(9, 16).Value = 0.005616140632441737
$ws.Cells.Item(9, 17).Value = 1.793119868727
$ws.Cells.Item(9, 18).Value = 16.138078818543
$ws.Cells.Item(9, 19).Value = 0.001133092304806193
$ws.Cells.Item(9, 20).Value = 0.001133092304806193
$ws.Cells.Item(10, 7).Value = 82.28866066666666
$ws.Cells.Item(10, 8).Value = 246.865982
$ws.Cells.Item(10, 9).Value = 0.177853633995942
$ws.Cells.Item(10, 10).Value = 0.177853633995942
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 0.8109183333333333
$ws.Cells.Item(10, 14).Value = 2.432755
$ws.Cells.Item(10, 15).Value = 0.2370884169621149
$ws.Cells.Item(10, 16).Value = 0.2370884169621149
$ws.Cells.Item(10, 17).Value = 66.72938356004555
$ws.Cells.Item(10, 18).Value = 600.5644520404098
$ws.Cells.Item(10, 19).Value = 0.04216703653505727
$ws.Cells.Item(10, 20).Value = 0.04216703653505727
$ws.Cells.Item(11, 7).Value = 82.28866066666666
$ws.Cells.Item(11, 8).Value = 246.865982
$ws.Cells.Item(11, 9).Value = 0.177853633995942
$ws.Cells.Item(11, 10).Value = 0.177853633995942
$ws.Cells.Item(11, 15).Value = 0.3378801459239538
$ws.Cells.Item(11, 16).Value = 0.3378801459239539
$ws.Cells.Item(11, 17).Value = 95.0975764382722
$ws.Cells.Item(11, 18).Value = 855.8781879444498
$ws.Cells.Item(11, 19).Value = 0.06009321180765435
$ws.Cells.Item(11, 20).Value = 0.06009321180765436
$ws.Cells.Item(12, 7).Value = 82.28866066666666
$ws.Cells.Item(12, 8).Value = 246.865982
$ws.Cells.Item(12, 9).Value = 0.177853633995942
$ws.Cells.Item(12, 10).Value = 0.177853633995942
$ws.Cells.Item(12, 13).Value = 1.434534666666667
$ws.Cells.Item(12, 14).Value = 4.303604
$ws.Cells.Item(12, 15).Value = 0.4194152964814894
$ws.Cells.Item(12, 16).Value = 0.4194152964814894
$ws.Cells.Item(12, 17).Value = 118.0459363999031
$ws.Cells.Item(12, 18).Value = 1062.413427599128
$ws.Cells.Item(12, 19).Value = 0.07459453463271831
$ws.Cells.Item(12, 20).Value = 0.07459453463271831
$ws.Cells.Item(13, 7).Value = 82.28866066666666
$ws.Cells.Item(13, 8).Value = 246.865982
$ws.Cells.Item(13, 9).Value = 0.177853633995942
$ws.Cells.Item(13, 10).Value = 0.177853633995942
$ws.Cells.Item(13, 11).Value = 1
$ws.Cells.Item(13, 12).Value = 0.3333333333333333
$ws.Cells.Item(13, 13).Value = 0.019209
$ws.Cells.Item(13, 14).Value = 0.057627
$ws.Cells.Item(13, 15).Value = 0.005616140632441737
$ws.Cells.Item(13, 16).Value = 0.005616140632441737
$ws.Cells.Item(13, 17).Value = 1.580682882746
$ws.Cells.Item(13, 18).Value = 14.226145944714
$ws.Cells.Item(13, 19).Value = 0.0009988510205120307
$ws.Cells.Item(13, 20).Value = 0.0009988510205120309
$ws.Cells.Item(14, 7).Value = 45.42364
$ws.Cells.Item(14, 8).Value = 136.27092
$ws.Cells.Item(14, 9).Value = 0.09817585288024938
$ws.Cells.Item(14, 10).Value = 0.09817585288024938
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 12).Value = 1
$ws.Cells.Item(14, 13).Value = 0.8109183333333333
$ws.Cells.Item(14, 14).Value = 2.432755
$ws.Cells.Item(14, 15).Value = 0.2370884169621149
$ws.Cells.Item(14, 16).Value = 0.2370884169621149
$ws.Cells.Item(14, 17).Value = 36.83486244273333
$ws.Cells.Item(14, 18).Value = 331.5137619846
$ws.Cells.Item(14, 19).Value = 0.02327635754328381
$ws.Cells.Item(14, 20).Value = 0.02327635754328382
$ws.Cells.Item(15, 7).Value = 45.42364
$ws.Cells.Item(15, 8).Value = 136.27092
$ws.Cells.Item(15, 9).Value = 0.09817585288024938
$ws.Cells.Item(15, 10).Value = 0.09817585288024938
$ws.Cells.Item(15, 15).Value = 0.3378801459239538
$ws.Cells.Item(15, 16).Value = 0.3378801459239539
$ws.Cells.Item(15, 17).Value = 52.49420809633332
$ws.Cells.Item(15, 18).Value = 472.4478728669999
$ws.Cells.Item(15, 19).Value = 0.03317167149738728
$ws.Cells.Item(15, 20).Value = 0.03317167149738729
$ws.Cells.Item(16, 7).Value = 45.42364
$ws.Cells.Item(16, 8).Value = 136.27092
$ws.Cells.Item(16, 9).Value = 0.09817585288024938
$ws.Cells.Item(16, 10).Value = 0.09817585288024938
$ws.Cells.Item(16, 13).Value = 1.434534666666667
$ws.Cells.Item(16, 14).Value = 4.303604
$ws.Cells.Item(16, 15).Value = 0.4194152964814894
$ws.Cells.Item(16, 16).Value = 0.4194152964814894
$ws.Cells.Item(16, 17).Value = 65.16178626618667
$ws.Cells.Item(16, 18).Value = 586.4560763956799
$ws.Cells.Item(16, 19).Value = 0.04117645444309288
$ws.Cells.Item(16, 20).Value = 0.04117645444309288
$ws.Cells.Item(17, 7).Value = 45.42364
$ws.Cells.Item(17, 8).Value = 136.27092
$ws.Cells.Item(17, 9).Value = 0.09817585288024938
$ws.Cells.Item(17, 10).Value = 0.09817585288024938
$ws.Cells.Item(17, 11).Value = 1
$ws.Cells.Item(17, 12).Value = 0.3333333333333333
$ws.Cells.Item(17, 13).Value = 0.019209
$ws.Cells.Item(17, 14).Value = 0.057627
$ws.Cells.Item(17, 15).Value = 0.005616140632441737
$ws.Cells.Item(17, 16).Value = 0.005616140632441737
$ws.Cells.Item(17, 17).Value = 0.87254270076
$ws.Cells.Item(17, 18).Value = 7.852884306839999
$ws.Cells.Item(17, 19).Value = 0.0005513693964853906
$ws.Cells.Item(17, 20).Value = 0.0005513693964853908
$ws.Cells.Item(18, 7).Value = 63.60851399999999
$ws.Cells.Item(18, 8).Value = 190.825542
$ws.Cells.Item(18, 9).Value = 0.1374795175462663
$ws.Cells.Item(18, 10).Value = 0.1374795175462663
$ws.Cells.Item(18, 11).Value = 3
$ws.Cells.Item(18, 12).Value = 1
$ws.Cells.Item(18, 13).Value = 0.8109183333333333
$ws.Cells.Item(18, 14).Value = 2.432755
$ws.Cells.Item(18, 15).Value = 0.2370884169621149
$ws.Cells.Item(18, 16).Value = 0.2370884169621149
$ws.Cells.Item(18, 17).Value = 51.58131015868999
$ws.Cells.Item(18, 18).Value = 464.2317914282099
$ws.Cells.Item(18, 19).Value = 0.03259480117975957
$ws.Cells.Item(18, 20).Value = 0.03259480117975957
$ws.Cells.Item(19, 7).Value = 63.60851399999999
$ws.Cells.Item(19, 8).Value = 190.825542
$ws.Cells.Item(19, 9).Value = 0.1374795175462663
$ws.Cells.Item(19, 10).Value = 0.1374795175462663
$ws.Cells.Item(19, 15).Value = 0.3378801459239538
$ws.Cells.Item(19, 16).Value = 0.3378801459239539
$ws.Cells.Item(19, 17).Value = 73.50970927504999
$ws.Cells.Item(19, 18).Value = 661.5873834754499
$ws.Cells.Item(19, 19).Value = 0.04645159945008722
$ws.Cells.Item(19, 20).Value = 0.04645159945008723
$ws.Cells.Item(20, 7).Value = 63.60851399999999
$ws.Cells.Item(20, 8).Value = 190.825542
$ws.Cells.Item(20, 9).Value = 0.1374795175462663
$ws.Cells.Item(20, 10).Value = 0.1374795175462663
$ws.Cells.Item(20, 13).Value = 1.434534666666667
$ws.Cells.Item(20, 14).Value = 4.303604
$ws.Cells.Item(20, 15).Value = 0.4194152964814894
$ws.Cells.Item(20, 16).Value = 0.4194152964814894
$ws.Cells.Item(20, 17).Value = 91.248618428152
$ws.Cells.Item(20, 18).Value = 821.237565853368
$ws.Cells.Item(20, 19).Value = 0.05766101261179939
$ws.Cells.Item(20, 20).Value = 0.05766101261179939
$ws.Cells.Item(21, 7).Value = 63.60851399999999
$ws.Cells.Item(21, 8).Value = 190.825542
$ws.Cells.Item(21, 9).Value = 0.1374795175462663
$ws.Cells.Item(21, 10).Value = 0.1374795175462663
$ws.Cells.Item(21, 11).Value = 1
$ws.Cells.Item(21, 12).Value = 0.3333333333333333
$ws.Cells.Item(21, 13).Value = 0.019209
$ws.Cells.Item(21, 14).Value = 0.057627
$ws.Cells.Item(21, 15).Value = 0.005616140632441737
$ws.Cells.Item(21, 16).Value = 0.005616140632441737
$ws.Cells.Item(21, 17).Value = 1.221855945426
$ws.Cells.Item(21, 18).Value = 10.996703508834
$ws.Cells.Item(21, 19).Value = 0.0007721043046200727
$ws.Cells.Item(21, 20).Value = 0.0007721043046200728
$ws.Cells.Item(22, 7).Value = 97.40706899999999
$ws.Cells.Item(22, 8).Value = 292.221207
$ws.Cells.Item(22, 9).Value = 0.2105296289694155
$ws.Cells.Item(22, 10).Value = 0.2105296289694155
$ws.Cells.Item(22, 11).Value = 3
$ws.Cells.Item(22, 12).Value = 1
$ws.Cells.Item(22, 13).Value = 0.8109183333333333
$ws.Cells.Item(22, 14).Value = 2.432755
$ws.Cells.Item(22, 15).Value = 0.2370884169621149
$ws.Cells.Item(22, 16).Value = 0.2370884169621149
$ws.Cells.Item(22, 17).Value = 78.98917804836499
$ws.Cells.Item(22, 18).Value = 710.902602435285
$ws.Cells.Item(22, 19).Value = 0.04991413645598013
$ws.Cells.Item(22, 20).Value = 0.04991413645598013
$ws.Cells.Item(23, 7).Value = 97.40706899999999
$ws.Cells.Item(23, 8).Value = 292.221207
$ws.Cells.Item(23, 9).Value = 0.2105296289694155
$ws.Cells.Item(23, 10).Value = 0.2105296289694155
$ws.Cells.Item(23, 15).Value = 0.3378801459239538
$ws.Cells.Item(23, 16).Value = 0.3378801459239539
$ws.Cells.Item(23, 17).Value = 112.569291015425
$ws.Cells.Item(23, 18).Value = 1013.123619138825
$ws.Cells.Item(23, 19).Value = 0.07113378175750196
$ws.Cells.Item(23, 20).Value = 0.07113378175750197
$ws.Cells.Item(24, 7).Value = 97.40706899999999
$ws.Cells.Item(24, 8).Value = 292.221207
$ws.Cells.Item(24, 9).Value = 0.2105296289694155
$ws.Cells.Item(24, 10).Value = 0.2105296289694155
$ws.Cells.Item(24, 13).Value = 1.434534666666667
$ws.Cells.Item(24, 14).Value = 4.303604
$ws.Cells.Item(24, 15).Value = 0.4194152964814894
$ws.Cells.Item(24, 16).Value = 0.4194152964814894
$ws.Cells.Item(24, 17).Value = 139.733817258892
$ws.Cells.Item(24, 18).Value = 1257.604355330028
$ws.Cells.Item(24, 19).Value = 0.08829934675234535
$ws.Cells.Item(24, 20).Value = 0.08829934675234535
$ws.Cells.Item(25, 7).Value = 97.40706899999999
$ws.Cells.Item(25, 8).Value = 292.221207
$ws.Cells.Item(25, 9).Value = 0.2105296289694155
$ws.Cells.Item(25, 10).Value = 0.2105296289694155
$ws.Cells.Item(25, 11).Value = 1
$ws.Cells.Item(25, 12).Value = 0.3333333333333333
$ws.Cells.Item(25, 13).Value = 0.019209
$ws.Cells.Item(25, 14).Value = 0.057627
$ws.Cells.Item(25, 15).Value = 0.005616140632441737
$ws.Cells.Item(25, 16).Value = 0.005616140632441737
$ws.Cells.Item(25, 17).Value = 1.871092388421
$ws.Cells.Item(25, 18).Value = 16.839831495789
$ws.Cells.Item(25, 19).Value = 0.001182364003588017
$ws.Cells.Item(25, 20).Value = 0.001182364003588017
